$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (D) and "Volume(1h)" (E) columns with refreshed crypto quote data.
# Each row: row number, new Price text, new Volume(1h) text.
# Price values that are plain decimals (e.g. "602.90") are written with a leading
# apostrophe so Excel stores them as text (matching the source file's inlineStr cells
# for the Price column) instead of silently converting them to numbers, and the cell
# style is reset to Normal afterwards so no stray quote-prefix formatting is left behind.
$rows = @(
    ,@(2, '65.627.07', '  +0.75%  ')
    ,@(3, '3.584.50', '  +1.50%  ')
    ,@(4, '1.00', '  -0.02%  ')
    ,@(5, '602.90', '  +0.38%  ')
    ,@(6, '136.71', '  -1.54%  ')
    ,@(7, '3.584.04', '  +1.39%  ')
    ,@(8, '1.00', '  +0.03%  ')
    ,@(9, '0.498', '  +1.66%  ')
    ,@(10, '0.125', '  +0.75%  ')
    ,@(11, '7.24', '  +4.84%  ')
    ,@(12, '0.391', '  +0.22%  ')
    ,@(13, '4.206.10', '  +1.75%  ')
    ,@(14, '28.01', '  +3.32%  ')
    ,@(15, '0.0000186', '  +0.92%  ')
    ,@(16, '3.593.87', '  +1.53%  ')
    ,@(17, '0.117', '  -0.15%  ')
    ,@(18, '65.763.18', '  +0.77%  ')
    ,@(19, '10.02', '  -2.72%  ')
    ,@(20, '14.60', '  +2.15%  ')
    ,@(21, '5.89', '  -0.79%  ')
    ,@(22, '395.55', '  +0.82%  ')
    ,@(23, '0.589', '  +2.99%  ')
    ,@(24, '3.734.14', '  +1.65%  ')
    ,@(25, '74.27', '  +0.70%  ')
    ,@(26, '1.00', '  +0.07%  ')
    ,@(27, '0.0000118', '  +3.00%  ')
    ,@(28, '8.17', '  +6.24%  ')
    ,@(29, '1.66', '  +30.78%  ')
    ,@(30, '2.41', '  +5.37%  ')
    ,@(31, '8.60', '  +4.83%  ')
    ,@(32, '1.00', '  +0.10%  ')
    ,@(33, '3.589.26', '  +1.25%  ')
    ,@(34, '24.49', '  +3.17%  ')
    ,@(35, '1.00', '  +0.01%  ')
    ,@(36, '0.147', '  +1.46%  ')
    ,@(37, '5.40', '  +8.80%  ')
    ,@(38, '1.61', '  +3.70%  ')
    ,@(39, '7.08', '  +1.84%  ')
    ,@(40, '168.29', '  -0.16%  ')
    ,@(41, '0.0836', '  +4.24%  ')
    ,@(42, '0.840', '  +2.07%  ')
    ,@(43, '26.59', '  +1.64%  ')
    ,@(44, '1.27', '  +7.91%  ')
    ,@(45, '43.26', '  +1.14%  ')
    ,@(46, '4.54', '  +2.65%  ')
    ,@(47, '1.00', '  -0.02%  ')
    ,@(48, '1.69', '  +1.40%  ')
    ,@(49, '7.04', '  +3.62%  ')
    ,@(50, '2.441.79', '  +1.22%  ')
    ,@(51, '316.54', '  +5.02%  ')
)

foreach ($row in $rows) {
    $r = $row[0]
    $price = $row[1]
    $volume = $row[2]

    if ($price -match "^[+-]?\d+(\.\d+)?$") {
        # Plain-looking decimal number -> force text entry via apostrophe prefix,
        # then restore the default (unstyled) cell style.
        $ws.Range("D" + $r).Value = "'" + $price
        $ws.Range("D" + $r).Style = "Normal"
    } else {
        $ws.Range("D" + $r).Value = $price
    }

    $ws.Range("E" + $r).Value = $volume
}
